# Update the "取得日時" (acquisition timestamp) column for all data rows
# on the "ランサーズ" sheet from 2025-12-04 06:28:42 to 2025-12-04 06:37:29.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2025-12-04 06:28:42"
$newValue = "2025-12-04 06:37:29"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
